$d = $word.ActiveDocument

# Start from the last paragraph in the document (the "All the time." one)
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range

# Insert the first new paragraph: "It is me again"
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.InsertAfter("It is me again")

# Insert the second new paragraph: "Wie geht’s."
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.InsertAfter("Wie geht" + [char]0x2019 + "s.")
